$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
  0.15885515775920567,
  -0.0059999999669528847,
  -0.0039999999732707181,
  -0.0079999999506661368,
  -0.0029999999759615648,
  0.048840816333273196,
  -0.0099999999362325731,
  -0.0099999999382598403,
  -0.0019999999847661876,
  0.057421440591115314,
  -0.0029999999836514135,
  -0.0034999999807077131,
  -0.003499999979666768,
  -0.0079999999555617762,
  -0.00099999999292865738,
  -0.0019999999878912433,
  -0.0019999999887314601,
  0.0086795838954012439,
  -0.047744225891290792,
  -0.0039999999750719439,
  -0.0039999999747948323,
  -0.0039999999746331838,
  -0.0049999999661878292,
  -0.01999999988312684,
  -0.019999999881517461,
  -0.0024999999739829804,
  -0.0024999999713664067,
  -0.0019999999638438126,
  -0.0069999999293299808,
  -0.059999999643558954,
  -0.0069999999228098631,
  -0.0099999999064976919,
  -0.0039999999383049101
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
